# Insert a new data row at row 31 (pushing existing rows 31..92 down to 32..93)
# and populate it with the new weekly price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Range("A31").Value2 = 11
$ws.Range("B31").Value2 = "Vega Monumental Concepción"
$ws.Range("C31").Value2 = "Bíobío"
$ws.Range("D31").Value2 = 44720
$ws.Range("E31").Value2 = 8
$ws.Range("F31").Value2 = 100112001
$ws.Range("G31").Value2 = "Berenjena"
$ws.Range("H31").Value2 = "Sin especificar"
$ws.Range("I31").Value2 = "Primera"
$ws.Range("J31").Value2 = 110
$ws.Range("K31").Value2 = 7500
$ws.Range("L31").Value2 = 8000
$ws.Range("M31").Value2 = 7727
$ws.Range("N31").Value2 = "`$/caja 60 unidades"
$ws.Range("O31").Value2 = "Región de Arica y Parinacota"
$ws.Range("P31").Value2 = 129
$ws.Range("Q31").Value2 = 60
$ws.Range("R31").Value2 = "Hortaliza"

$ws.Range("D31").NumberFormat = $ws.Range("D32").NumberFormat
